# Khulo.xlsx - "upgrade left table until javakheti"
#
# This sheet had a placeholder tab name ("1") and a left-hand data table
# (Total / Urban / Rural rows for 2010-2023) in which several "Urban" and
# "Rural" counts were actually unavailable/confidential. Bring the sheet
# in line with the rest of the workbook series: rename the tab, mark the
# unavailable Urban/Rural cells with the confidential-data marker used
# elsewhere in the table, and remove the blank spacer row above the note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the worksheet tab its proper name instead of the default "1"
$ws.Name = "Khulo"

$ellipsis = [char]0x2026

# Row 6 ("Urban"): almost every year becomes confidential/unavailable.
# B6 uses a plain three-dot marker (new distinct shared string), the rest
# use the existing "…" ellipsis marker already used throughout the sheet.
$ws.Range("B6").Value2 = "..."
$ws.Range("C6").Value2 = $ellipsis
$ws.Range("D6").Value2 = $ellipsis
$ws.Range("F6").Value2 = $ellipsis
$ws.Range("H6").Value2 = $ellipsis
$ws.Range("I6").Value2 = $ellipsis
$ws.Range("K6").Value2 = $ellipsis
$ws.Range("M6").Value2 = $ellipsis

# Row 7 ("Rural"): 2016, 2017, 2019 and 2021 become unavailable too.
$ws.Range("H7").Value2 = $ellipsis
$ws.Range("I7").Value2 = $ellipsis
$ws.Range("K7").Value2 = $ellipsis
$ws.Range("M7").Value2 = $ellipsis

# Remove the now-unneeded blank spacer row between the data table and the
# footnote, shifting the note up from row 9 to row 8.
$ws.Rows("8:8").Delete()
